$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 - column headers (bold "title" style, same as A3)
$ws.Cells.Item(9, 2).Value = "Number of employees"
$ws.Cells.Item(9, 2).Font.Bold = $true

$ws.Cells.Item(9, 3).Value = "Assets (local currency, unless noted otherwise)"
$ws.Cells.Item(9, 3).Font.Bold = $true

$ws.Cells.Item(9, 4).Value = "Turnover (local currency, unless noted otherwise)"
$ws.Cells.Item(9, 4).Font.Bold = $true

# Row 10 - Micro
$ws.Cells.Item(10, 1).Value = "Micro"
$ws.Cells.Item(10, 2).Value = "0-3"
$ws.Cells.Item(10, 3).Value = "'"
$ws.Cells.Item(10, 3).Style = "Normal"
$ws.Cells.Item(10, 4).Value = "'"
$ws.Cells.Item(10, 4).Style = "Normal"

# Row 11 - Small
$ws.Cells.Item(11, 1).Value = "Small"
$ws.Cells.Item(11, 2).Value = "4-20"
$ws.Cells.Item(11, 3).Value = "'"
$ws.Cells.Item(11, 3).Style = "Normal"
$ws.Cells.Item(11, 4).Value = "'"
$ws.Cells.Item(11, 4).Style = "Normal"

# Row 12 - Medium
$ws.Cells.Item(12, 1).Value = "Medium"
$ws.Cells.Item(12, 2).Value = "21-50"
$ws.Cells.Item(12, 3).Value = "'"
$ws.Cells.Item(12, 3).Style = "Normal"
$ws.Cells.Item(12, 4).Value = "'"
$ws.Cells.Item(12, 4).Style = "Normal"

# Row 13 - Large
$ws.Cells.Item(13, 1).Value = "Large"
$ws.Cells.Item(13, 2).Value = ">50"
$ws.Cells.Item(13, 3).Value = "'"
$ws.Cells.Item(13, 3).Style = "Normal"
$ws.Cells.Item(13, 4).Value = "'"
$ws.Cells.Item(13, 4).Style = "Normal"
